# Auto-generated edit script: update crypto price/volume values per diff
# (generated from canonical-OOXML before/after comparison)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cells whose new value is plain text / not a valid Excel number: set directly ---
$ws.Range('D2').Value = '26.579.07'
$ws.Range('E2').Value = '  +0.38%  '
$ws.Range('D3').Value = '1.740.05'
$ws.Range('E3').Value = '  +0.58%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('E5').Value = '  +1.29%  '
$ws.Range('E7').Value = '  +2.69%  '
$ws.Range('E8').Value = '  +0.53%  '
$ws.Range('E9').Value = '  +1.02%  '
$ws.Range('D10').Value = '1.736.20'
$ws.Range('E11').Value = '  -0.82%  '
$ws.Range('E12').Value = '  +0.10%  '
$ws.Range('E13').Value = '  -0.45%  '
$ws.Range('E14').Value = '  +1.00%  '
$ws.Range('E15').Value = '  +1.20%  '
$ws.Range('E16').Value = '  -0.02%  '
$ws.Range('B17').Value = 'ShibaInu'
$ws.Range('C17').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('E17').Value = '  +6.29%  '
$ws.Range('B18').Value = 'WrappedBTC'
$ws.Range('C18').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D18').Value = '26.585.04'
$ws.Range('E18').Value = '  +0.37%  '
$ws.Range('E19').Value = '  +0.04%  '
$ws.Range('D21').Value = '1.959.24'
$ws.Range('E21').Value = '  +0.29%  '
$ws.Range('E22').Value = '  +0.86%  '
$ws.Range('E23').Value = '  -1.94%  '
$ws.Range('E24').Value = '  -1.03%  '
$ws.Range('E25').Value = '  +2.45%  '
$ws.Range('E26').Value = '  +0.90%  '
$ws.Range('E27').Value = '  +0.66%  '
$ws.Range('E28').Value = '  -1.44%  '
$ws.Range('E29').Value = '  +1.13%  '
$ws.Range('E30').Value = '  +1.57%  '
$ws.Range('E31').Value = '  +0.97%  '
$ws.Range('E32').Value = '  -0.40%  '
$ws.Range('E33').Value = '  +1.30%  '
$ws.Range('E34').Value = '  -0.03%  '
$ws.Range('E35').Value = '  -0.07%  '
$ws.Range('E36').Value = '  +2.88%  '
$ws.Range('E37').Value = '  -0.43%  '
$ws.Range('E38').Value = '  -4.66%  '
$ws.Range('E39').Value = '  +2.66%  '
$ws.Range('E40').Value = '  -0.06%  '
$ws.Range('E41').Value = '  -0.10%  '
$ws.Range('E42').Value = '  +0.18%  '
$ws.Range('E43').Value = '  -5.06%  '
$ws.Range('E44').Value = '  -4.22%  '
$ws.Range('E45').Value = '  +0.60%  '
$ws.Range('E46').Value = '  -0.29%  '
$ws.Range('E47').Value = '  -0.34%  '
$ws.Range('E48').Value = '  +1.32%  '
$ws.Range('E49').Value = '  -0.57%  '
$ws.Range('E50').Value = '  -0.65%  '
$ws.Range('E51').Value = '  +0.04%  '

# --- Cells whose new value LOOKS like a number (e.g. "247.09"): must stay text ---
# Force text format first so Excel keeps the literal string instead of converting
# it to a numeric value (which would also round/alter the trailing zeros), then
# reset the cell style back to Normal so no stray number-format style lingers.
$numericLookingCells = @('D4', 'D5', 'D6', 'D7', 'D8', 'D11', 'D12', 'D13', 'D14', 'D15', 'D16', 'D17', 'D19', 'D22', 'D23', 'D24', 'D25', 'D26', 'D28', 'D29', 'D30', 'D31', 'D32', 'D33', 'D34', 'D36', 'D37', 'D38', 'D39', 'D40', 'D43', 'D44', 'D45', 'D46', 'D47', 'D49', 'D50', 'D51')
foreach ($addr in $numericLookingCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D4').Value = '0.9997'
$ws.Range('D5').Value = '247.09'
$ws.Range('D6').Value = '0.9997'
$ws.Range('D7').Value = '0.4930'
$ws.Range('D8').Value = '0.2685'
$ws.Range('D11').Value = '0.07055'
$ws.Range('D12').Value = '15.75'
$ws.Range('D13').Value = '0.6153'
$ws.Range('D14').Value = '4.590'
$ws.Range('D15').Value = '78.01'
$ws.Range('D16').Value = '0.9996'
$ws.Range('D17').Value = '0.000007363'
$ws.Range('D19').Value = '0.9999'
$ws.Range('D22').Value = '4.604'
$ws.Range('D23').Value = '8.739'
$ws.Range('D24').Value = '5.261'
$ws.Range('D25').Value = '139.74'
$ws.Range('D26').Value = '15.48'
$ws.Range('D28').Value = '1.765'
$ws.Range('D29').Value = '107.86'
$ws.Range('D30').Value = '4.052'
$ws.Range('D31').Value = '0.08063'
$ws.Range('D32').Value = '3.730'
$ws.Range('D33').Value = '0.04623'
$ws.Range('D34').Value = '0.9992'
$ws.Range('D36').Value = '1.019'
$ws.Range('D37').Value = '0.6383'
$ws.Range('D38').Value = '0.8998'
$ws.Range('D39').Value = '2.042'
$ws.Range('D40').Value = '2.405'
$ws.Range('D43').Value = '101.93'
$ws.Range('D44').Value = '5.414'
$ws.Range('D45').Value = '0.3933'
$ws.Range('D46').Value = '6.909'
$ws.Range('D47').Value = '0.1188'
$ws.Range('D49').Value = '7.856'
$ws.Range('D50').Value = '30.61'
$ws.Range('D51').Value = '1.272'

foreach ($addr in $numericLookingCells) {
    $ws.Range($addr).Style = "Normal"
}
